$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Kostenstelle) used to hold text codes like "TS-001"; now holds
# plain numeric cost-center IDs (1001-1020) instead.
$ws.Range("D2").Value  = 1001
$ws.Range("D3").Value  = 1002
$ws.Range("D4").Value  = 1003
$ws.Range("D5").Value  = 1004
$ws.Range("D6").Value  = 1005
$ws.Range("D7").Value  = 1006
$ws.Range("D8").Value  = 1007
$ws.Range("D9").Value  = 1008
$ws.Range("D10").Value = 1009
$ws.Range("D11").Value = 1010
$ws.Range("D12").Value = 1011
$ws.Range("D13").Value = 1012
$ws.Range("D14").Value = 1013
$ws.Range("D15").Value = 1014
$ws.Range("D16").Value = 1015
$ws.Range("D17").Value = 1016
$ws.Range("D18").Value = 1017
$ws.Range("D19").Value = 1018
$ws.Range("D20").Value = 1019
$ws.Range("D21").Value = 1020

# With the long text codes gone, every data row now wraps to the same,
# shorter height.
$ws.Rows("2:21").RowHeight = 29.25

# Page now prints on A4 portrait paper.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moved onto the (now numeric) Kostenstelle column.
$ws.Range("D2:D21").Select() | Out-Null
